$d = $word.ActiveDocument

# Locate the paragraph that ends the document: "When you look to the staaaars"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*staaaars*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph containing 'staaaars'"
}

$full = $target.Range
# Exclude the trailing paragraph mark from the replaced range so that
# InsertXML doesn't leave behind a stray empty paragraph (the body's
# very last paragraph mark is structurally "sticky").
$r = $d.Range($full.Start, $full.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">When you look to the </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>staaaars</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
<w:r><w:lastRenderedPageBreak/><w:t>As much as we love these big attitude Texans, they had it coming</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$r.InsertXML($xml)
